# Insert a new weekly price-report row for "Chirimoya" at row 205,
# pushing the existing rows 205:215 down to 206:216.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 205:215 down to 206:216 by inserting a new blank row at 205.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new record.
$ws.Range("A205").Value = 5
$ws.Range("B205").Value = "Macroferia Regional de Talca"
$ws.Range("C205").Value = "Maule"
$ws.Range("D205").Value = (Get-Date -Year 2023 -Month 12 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D205").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E205").Value = 7
$ws.Range("F205").Value = "Fruta"
$ws.Range("G205").Value = 100107
$ws.Range("H205").Value = "Otros"
$ws.Range("I205").Value = 100107002
$ws.Range("J205").Value = "Chirimoya"
$ws.Range("K205").Value = "Cultivar IV Región"
$ws.Range("L205").Value = "Primera"
$ws.Range("M205").Value = 180
$ws.Range("N205").Value = 17000
$ws.Range("O205").Value = 17000
$ws.Range("P205").Value = 17000
$ws.Range("Q205").Value = "$/bandeja 10 kilos"
$ws.Range("R205").Value = "Provincia de Limarí"
$ws.Range("S205").Value = 1700
$ws.Range("T205").Value = 10
